{"js": "// Add an \"Ethical Considerations\" section to the end of the write-up.\n// The document currently ends with a trailing empty paragraph right\n// before the section break; we turn that spot into two new paragraphs:\n//   1) a bold, 18pt \"Ethical Considerations\" heading (matching the\n//      other section headings in the doc), and\n//   2) a body paragraph (Segoe UI, color #24292E, white shading - the\n//      formatting that came along when the text was pasted in from a\n//      web source) with the ethics prompt text.\n//\n// We build the two paragraphs as raw WordprocessingML and insert them\n// with Body.insertOoxml so the exact run-property markup (w:b, w:bCs,\n// w:sz/w:szCs, w:rFonts, w:shd, etc.) is reproduced faithfully.\n\nconst body = context.document.body;\n\nconst wAttr = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nconst headingParagraph =\n  `<w:p ${wAttr}>` +\n    `<w:r>` +\n      `<w:rPr>` +\n        `<w:b/><w:bCs/>` +\n        `<w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/>` +\n      `</w:rPr>` +\n      `<w:t>Ethical Considerations</w:t>` +\n    `</w:r>` +\n  `</w:p>`;\n\nconst bodyParagraph =\n  `<w:p ${wAttr}>` +\n    `<w:r>` +\n      `<w:rPr>` +\n        `<w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>` +\n        `<w:color w:val=\"24292E\"/>` +\n        `<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>` +\n      `</w:rPr>` +\n      `<w:t>I</w:t>` +\n    `</w:r>` +\n    `<w:r>` +\n      `<w:rPr>` +\n        `<w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>` +\n        `<w:color w:val=\"24292E\"/>` +\n        `<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>` +\n      `</w:rPr>` +\n      `<w:t>deas about the application of AI ethics that might be necessary to use your solution in real-world scenario</w:t>` +\n    `</w:r>` +\n  `</w:p>`;\n\nconst ooxmlPackage =\n  `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>` +\n  `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n    `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n      `<pkg:xmlData>` +\n        `<w:document ${wAttr}>` +\n          `<w:body>${headingParagraph}${bodyParagraph}</w:body>` +\n        `</w:document>` +\n      `</pkg:xmlData>` +\n    `</pkg:part>` +\n  `</pkg:package>`;\n\n// Inserting at the end of the body fills the existing trailing empty\n// paragraph with the first inserted paragraph's content, and appends\n// the rest, so we end up with exactly two new paragraphs (the stray\n// empty paragraph isn't left behind as a third one).\nbody.insertOoxml(ooxmlPackage, Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Add an \"Ethical Considerations\" section to the end of the write-up.\n#\n# The document currently ends with a trailing empty paragraph right\n# before the section break. We turn that spot into two new paragraphs:\n#   1) a bold, 18pt \"Ethical Considerations\" heading (matching the\n#      other section headings in the doc), and\n#   2) a body paragraph (Segoe UI, color #24292E, white \"clear\"\n#      shading - the formatting that tags along when text is pasted in\n#      from a web page) containing the ethics prompt text.\n#\n# We build the two paragraphs as raw WordprocessingML and insert them\n# with Range.InsertXML so the exact run-property markup (w:b, w:bCs,\n# w:sz/w:szCs, w:rFonts, w:shd, etc.) is reproduced faithfully.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$insertionRange = $lastParagraph.Range\n\n$wAttr = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n$headingParagraph = '<w:p ' + $wAttr + '>' + `\n    '<w:r>' + `\n      '<w:rPr>' + `\n        '<w:b/><w:bCs/>' + `\n        '<w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/>' + `\n      '</w:rPr>' + `\n      '<w:t>Ethical Considerations</w:t>' + `\n    '</w:r>' + `\n  '</w:p>'\n\n$bodyParagraph = '<w:p ' + $wAttr + '>' + `\n    '<w:r>' + `\n      '<w:rPr>' + `\n        '<w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' + `\n        '<w:color w:val=\"24292E\"/>' + `\n        '<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>' + `\n      '</w:rPr>' + `\n      '<w:t>I</w:t>' + `\n    '</w:r>' + `\n    '<w:r>' + `\n      '<w:rPr>' + `\n        '<w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' + `\n        '<w:color w:val=\"24292E\"/>' + `\n        '<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>' + `\n      '</w:rPr>' + `\n      '<w:t>deas about the application of AI ethics that might be necessary to use your solution in real-world scenario</w:t>' + `\n    '</w:r>' + `\n  '</w:p>'\n\n$ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n      '<pkg:xmlData>' + `\n        '<w:document ' + $wAttr + '>' + `\n          '<w:body>' + $headingParagraph + $bodyParagraph + '</w:body>' + `\n        '</w:document>' + `\n      '</pkg:xmlData>' + `\n    '</pkg:part>' + `\n  '</pkg:package>'\n\n# InsertXML inserts the two new paragraphs just before $insertionRange,\n# pushing the original (empty) trailing paragraph after them instead of\n# replacing it. Remove that now-redundant paragraph break by deleting\n# the paragraph mark at the end of the newly inserted body paragraph,\n# which merges it with the (still empty) old trailing paragraph.\n$insertionRange.InsertXML($ooxmlPackage)\n\n$mergeTarget = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n$markRange = $d.Range($mergeTarget.Range.End - 1, $mergeTarget.Range.End)\n$markRange.Delete()\n"}
